$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.037.26'
$ws.Range("E2").Value = '  +2.45%  '
$ws.Range("D3").Value = '2.280.80'
$ws.Range("E3").Value = '  +1.67%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '230.84'
$ws.Range("E5").Value = '  -0.41%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.636'
$ws.Range("E6").Value = '  +2.11%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '63.64'
$ws.Range("E7").Value = '  +4.42%  '
$ws.Range("E8").Value = '  +0.09%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.447'
$ws.Range("E9").Value = '  +9.20%  '
$ws.Range("E10").Value = '  +10.09%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '56.99'
$ws.Range("E11").Value = '  -1.51%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '26.55'
$ws.Range("E12").Value = '  +18.16%  '
$ws.Range("E13").Value = '  +2.05%  '
$ws.Range("D14").Value = '2.614.65'
$ws.Range("E14").Value = '  +1.57%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.68'
$ws.Range("E15").Value = '  +0.49%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.17'
$ws.Range("E16").Value = '  +8.40%  '
$ws.Range("E17").Value = '  +4.93%  '
$ws.Range("D18").Value = '2.275.25'
$ws.Range("E18").Value = '  +1.58%  '
$ws.Range("D19").Value = '43.888.49'
$ws.Range("E19").Value = '  +2.38%  '
$ws.Range("D20").Value = '0.0₃0998'
$ws.Range("E20").Value = '  +6.03%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '73.75'
$ws.Range("E21").Value = '  +1.51%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.09'
$ws.Range("E22").Value = '  -2.00%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '253.40'
$ws.Range("E23").Value = '  +2.97%  '
$ws.Range("E24").Value = '  +0.22%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.45'
$ws.Range("E25").Value = '  -5.85%  '
$ws.Range("E26").Value = '  -2.61%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.35'
$ws.Range("E27").Value = '  +25.75%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.03'
$ws.Range("E28").Value = '  +2.69%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '172.00'
$ws.Range("E29").Value = '  +1.31%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.139'
$ws.Range("E30").Value = '  -2.29%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.85'
$ws.Range("E31").Value = '  +2.08%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.39'
$ws.Range("E32").Value = '  -6.47%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.125'
$ws.Range("E33").Value = '  +3.26%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0698'
$ws.Range("E34").Value = '  +6.26%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.81'
$ws.Range("E35").Value = '  +1.70%  '
$ws.Range("E36").Value = '  -2.33%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.81'
$ws.Range("E37").Value = '  +5.58%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.56'
$ws.Range("E38").Value = '  +1.98%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.32'
$ws.Range("E39").Value = '  -2.73%  '
$ws.Range("E40").Value = '  +3.64%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.000239'
$ws.Range("E41").Value = '  +6.42%  '
$ws.Range("E42").Value = '  +0.18%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '17.76'
$ws.Range("E43").Value = '  +7.93%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0975'
$ws.Range("E44").Value = '  +0.47%  '
$ws.Range("E45").Value = '  -4.82%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '98.36'
$ws.Range("E46").Value = '  +1.17%  '
$ws.Range("B47").Value = 'Celestia'
$ws.Range("C47").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.33'
$ws.Range("E47").Value = '  +13.33%  '
$ws.Range("B48").Value = 'TrustWalletToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.20'
$ws.Range("E48").Value = '  -0.89%  '
$ws.Range("B49").Value = 'FTXToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.38'
$ws.Range("E49").Value = '  -0.29%  '
$ws.Range("D50").Value = '1.448.10'
$ws.Range("E50").Value = '  -1.13%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.29'
$ws.Range("E51").Value = '  +2.73%  '
